$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Apply the same fill/format used by the existing "zeroshot huang combined
# --- with own (w/o marketing)" rows (e.g. row 34) to the two new rows 37:38 ---
[void]$ws.Range("A34:L34").Copy()
[void]$ws.Range("A37:L38").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 37
$ws.Range("H37").Value = "5500s"
$ws.Range("A37").Value = "zeroshot huang combined with own (w/o marketing)"
$ws.Range("D37").Value = 2500
$ws.Range("E37").Value = 200
$ws.Range("F37").Value = 5
$ws.Range("G37").Value = 0.25
$ws.Range("I37").Value = "null"
$ws.Range("J37").Value = "yes"
$ws.Range("K37").Value = 77

# --- new "regular / full random / sentences" runs, rows 40-42 ---
[void]$ws.Range("D34").Copy()
[void]$ws.Range("D40").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("H40").Value = "4900s"
$ws.Range("A40").Value = "regular"
$ws.Range("B40").Value = "full random"
$ws.Range("C40").Value = "sentences"
$ws.Range("D40").Value = 2500
$ws.Range("E40").Value = 200
$ws.Range("F40").Value = 5
$ws.Range("I40").Value = 20
$ws.Range("J40").Value = "yes"
$ws.Range("K40").Value = 20

# Row 38 (second zeroshot run, filled in after the "regular" block)
$ws.Range("H38").Value = "5900s"
$ws.Range("A38").Value = "zeroshot huang combined with own (w/o marketing)"
$ws.Range("D38").Value = 2500
$ws.Range("E38").Value = 200
$ws.Range("F38").Value = 5
$ws.Range("G38").Value = 0.75
$ws.Range("I38").Value = "null"
$ws.Range("J38").Value = "yes"
$ws.Range("K38").Value = 144

$ws.Range("L40").Value = "medium (z.b. costs +earnings sind grouped)"

# Row 41
[void]$ws.Range("D34").Copy()
[void]$ws.Range("D41").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A41").Value = "regular"
$ws.Range("B41").Value = "full random"
$ws.Range("C41").Value = "sentences"
$ws.Range("D41").Value = 2500
$ws.Range("E41").Value = 200
$ws.Range("F41").Value = 5
$ws.Range("H41").Value = "4900s"
$ws.Range("I41").Value = 30
$ws.Range("J41").Value = "yes"
$ws.Range("K41").Value = 30

# Row 42
[void]$ws.Range("D34").Copy()
[void]$ws.Range("D42").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A42").Value = "regular"
$ws.Range("B42").Value = "full random"
$ws.Range("C42").Value = "sentences"
$ws.Range("D42").Value = 2500
$ws.Range("E42").Value = 200
$ws.Range("F42").Value = 5
$ws.Range("H42").Value = "4900s"
$ws.Range("I42").Value = 40
$ws.Range("J42").Value = "yes"
$ws.Range("K42").Value = 40
$ws.Range("L42").Value = "good"

# Match the saved selection / scroll position from the author's session.
[void]$excel.Goto($ws.Range("A17"), $true)
[void]$ws.Range("O37").Select()

"applied edit"
